# Update "想去人数" (F column) values across sheets, matching the
# refreshed data snapshot from the gh-pages generation run.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 976
$ws.Range("F4").Value = 220
$ws.Range("F5").Value = 25
$ws.Range("F6").Value = 1097
$ws.Range("F7").Value = 880
$ws.Range("F8").Value = 267
$ws.Range("F11").Value = 861
$ws.Range("F14").Value = 513
$ws.Range("F15").Value = 1355
$ws.Range("F16").Value = 122
$ws.Range("F17").Value = 1277
$ws.Range("F18").Value = 1227
$ws.Range("F19").Value = 2907
$ws.Range("F20").Value = 1491
$ws.Range("F21").Value = 733
$ws.Range("F23").Value = 1294
$ws.Range("F25").Value = 1048
$ws.Range("F26").Value = 366
$ws.Range("F27").Value = 3219
$ws.Range("F29").Value = 541
$ws.Range("F30").Value = 1434

# --- Sheet: 演出 ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 62

# --- Sheet: 本地生活 ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 760

# --- Sheet: 全部类型 ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 760
$ws.Range("F6").Value = 976
$ws.Range("F7").Value = 220
$ws.Range("F9").Value = 25
$ws.Range("F10").Value = 1097
$ws.Range("F11").Value = 880
$ws.Range("F12").Value = 267
$ws.Range("F13").Value = 62
$ws.Range("F21").Value = 861
$ws.Range("F24").Value = 513
$ws.Range("F25").Value = 1355
$ws.Range("F26").Value = 122
$ws.Range("F27").Value = 1277
$ws.Range("F28").Value = 1227
$ws.Range("F29").Value = 2907
$ws.Range("F30").Value = 1491
$ws.Range("F31").Value = 733
$ws.Range("F33").Value = 1294
$ws.Range("F37").Value = 1048
$ws.Range("F38").Value = 366
$ws.Range("F39").Value = 3219
$ws.Range("F41").Value = 541
$ws.Range("F42").Value = 1434
